$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (avoids Excel
# auto-converting numeric-looking strings like "1.003" or "9.080" into
# actual numbers, which would lose formatting / trailing zeros).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "29.222.19"
$ws.Range("E2").Value = "  +2.86%  "
$ws.Range("D3").Value = "1.901.17"
$ws.Range("E3").Value = "  +1.57%  "
Set-TextValue $ws.Range("D4") "1.003"
$ws.Range("E4").Value = "  -1.81%  "
Set-TextValue $ws.Range("D5") "315.45"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("E6").Value = "  -1.72%  "
Set-TextValue $ws.Range("D7") "0.5121"
$ws.Range("E7").Value = "  +0.33%  "
Set-TextValue $ws.Range("D8") "0.3941"
$ws.Range("E8").Value = "  -0.39%  "
Set-TextValue $ws.Range("D9") "0.08415"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D10") "42.59"
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D11") "1.122"
$ws.Range("E11").Value = "  +1.26%  "
Set-TextValue $ws.Range("D12") "6.262"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "1.895.92"
$ws.Range("E13").Value = "  +1.00%  "
Set-TextValue $ws.Range("D14") "20.54"
$ws.Range("E14").Value = "  +0.54%  "
Set-TextValue $ws.Range("D15") "7.359"
$ws.Range("E15").Value = "  +1.87%  "
Set-TextValue $ws.Range("D16") "1.003"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("E17").Value = "  +2.41%  "
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("E19").Value = "  -0.86%  "
Set-TextValue $ws.Range("D20") "17.87"
$ws.Range("E20").Value = "  +0.99%  "
Set-TextValue $ws.Range("D21") "1.004"
$ws.Range("E21").Value = "  -1.54%  "
Set-TextValue $ws.Range("D22") "6.011"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("D23").Value = "29.216.47"
$ws.Range("E23").Value = "  +2.57%  "
Set-TextValue $ws.Range("D24") "11.19"
$ws.Range("E24").Value = "  +0.33%  "
Set-TextValue $ws.Range("D25") "2.227"
$ws.Range("E25").Value = "  -2.87%  "
$ws.Range("D26").Value = "2.116.32"
$ws.Range("E26").Value = "  +1.32%  "
Set-TextValue $ws.Range("D27") "21.01"
$ws.Range("E27").Value = "  +1.44%  "
Set-TextValue $ws.Range("D28") "158.71"
$ws.Range("E28").Value = "  -1.93%  "
Set-TextValue $ws.Range("D29") "2.429"
$ws.Range("E29").Value = "  +3.87%  "
Set-TextValue $ws.Range("D30") "127.21"
$ws.Range("E30").Value = "  +0.27%  "
Set-TextValue $ws.Range("D31") "1.062"
$ws.Range("E31").Value = "  +2.59%  "
Set-TextValue $ws.Range("D32") "0.1047"
$ws.Range("E32").Value = "  -0.56%  "
Set-TextValue $ws.Range("D33") "5.913"
$ws.Range("E33").Value = "  +2.82%  "
Set-TextValue $ws.Range("D34") "3.631"
$ws.Range("E34").Value = "  +0.07%  "
Set-TextValue $ws.Range("D35") "0.02476"
$ws.Range("E35").Value = "  +1.78%  "
Set-TextValue $ws.Range("D36") "0.06614"
$ws.Range("E36").Value = "  +2.44%  "
Set-TextValue $ws.Range("D37") "9.080"
$ws.Range("E37").Value = "  +3.10%  "
Set-TextValue $ws.Range("D38") "0.2195"
$ws.Range("E38").Value = "  +0.88%  "
Set-TextValue $ws.Range("D39") "1.228"
$ws.Range("E39").Value = "  +4.04%  "
Set-TextValue $ws.Range("D40") "5.089"
$ws.Range("E40").Value = "  +2.31%  "
Set-TextValue $ws.Range("D41") "0.6469"
$ws.Range("E41").Value = "  +1.77%  "
Set-TextValue $ws.Range("D42") "1.237"
$ws.Range("E42").Value = "  -1.81%  "
Set-TextValue $ws.Range("D43") "11.29"
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("E44").Value = "  -1.61%  "
Set-TextValue $ws.Range("D45") "0.6038"
$ws.Range("E45").Value = "  +0.28%  "
Set-TextValue $ws.Range("D46") "13.19"
$ws.Range("E46").Value = "  +1.46%  "
Set-TextValue $ws.Range("D47") "3.676"
$ws.Range("E47").Value = "  -0.93%  "
Set-TextValue $ws.Range("D48") "2.045"
$ws.Range("E48").Value = "  +2.95%  "
$ws.Range("E49").Value = "  +2.43%  "
Set-TextValue $ws.Range("D50") "123.01"
$ws.Range("E50").Value = "  +1.02%  "
Set-TextValue $ws.Range("D51") "1.162"
$ws.Range("E51").Value = "  -3.71%  "
